# Update Betfair Back/Lay odds cells on the active sheet to match the
# latest data pull for 2025-11-19.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Brazilian Serie A: SE Palmeiras vs EC Vitoria Salvador
$ws.Range("F2").Value = 1.33
$ws.Range("G2").Value = 1.35
$ws.Range("J2").Value = 5.7
$ws.Range("K2").Value = 6
$ws.Range("N2").Value = 4.1
$ws.Range("O2").Value = 1.29
$ws.Range("P2").Value = 2.1
$ws.Range("Q2").Value = 1.83
$ws.Range("S2").Value = 3.15
$ws.Range("T2").Value = 2.3
$ws.Range("U2").Value = 1.69
$ws.Range("V2").Value = 1.08
$ws.Range("W2").Value = 3.85
$ws.Range("X2").Value = 19
$ws.Range("AB2").Value = 7.6
$ws.Range("AD2").Value = 46
$ws.Range("AF2").Value = 7.2
$ws.Range("AG2").Value = 10.5
$ws.Range("AI2").Value = 1000
$ws.Range("AL2").Value = 48
$ws.Range("AN2").Value = 6.4

# Row 3 - Brazilian Serie B: Vila Nova vs Volta Redonda
$ws.Range("F3").Value = 1.86
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 4.9

# Row 4 - Chilean Primera B: Santiago Wanderers vs Cobreloa Calama
$ws.Range("F4").Value = 2.46
$ws.Range("G4").Value = 2.92
$ws.Range("I4").Value = 3.5
$ws.Range("K4").Value = 3.7
$ws.Range("V4").Value = 1.44
$ws.Range("W4").Value = 1.52

# Row 5 - Colombian Primera A: Atletico Bucaramanga vs Santa Fe
$ws.Range("F5").Value = 1.95
$ws.Range("G5").Value = 2.14
$ws.Range("J5").Value = 2.84
$ws.Range("K5").Value = 3.45
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 2.68
$ws.Range("Q5").Value = 2.3
$ws.Range("W5").Value = 1.92
$ws.Range("X5").Value = 11
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 1000

# Row 6 - Brazilian Serie A: Fluminense vs Flamengo
$ws.Range("F6").Value = 4.1
$ws.Range("G6").Value = 4.6
$ws.Range("H6").Value = 2.08
$ws.Range("I6").Value = 2.18
$ws.Range("J6").Value = 3.2
$ws.Range("K6").Value = 3.4
$ws.Range("M6").Value = 1.11
$ws.Range("O6").Value = 1.5
$ws.Range("P6").Value = 1.61
$ws.Range("Q6").Value = 2.46
$ws.Range("R6").Value = 1.22
$ws.Range("S6").Value = 4.8
$ws.Range("T6").Value = 2.1
$ws.Range("U6").Value = 1.79
$ws.Range("V6").Value = 1.86
$ws.Range("W6").Value = 1.28

# Row 7 - Brazilian Serie A: Santos vs Mirassol
$ws.Range("G7").Value = 2.26
$ws.Range("H7").Value = 3.55
$ws.Range("I7").Value = 3.85
$ws.Range("J7").Value = 3.45
$ws.Range("K7").Value = 3.7
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 1.84
$ws.Range("R7").Value = 1.32
$ws.Range("V7").Value = 1.35
$ws.Range("W7").Value = 1.79
$ws.Range("Y7").Value = 13
$ws.Range("AD7").Value = 980

# Row 8 - Brazilian Serie A: Gremio vs Vasco da Gama
$ws.Range("G8").Value = 2.36
$ws.Range("K8").Value = 3.6
$ws.Range("M8").Value = 1.08
$ws.Range("N8").Value = 3.3
$ws.Range("O8").Value = 1.38
$ws.Range("P8").Value = 1.79
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("T8").Value = 1.86
$ws.Range("U8").Value = 2.02
$ws.Range("V8").Value = 1.37
$ws.Range("W8").Value = 1.73
$ws.Range("X8").Value = 980
$ws.Range("Y8").Value = 12.5
$ws.Range("Z8").Value = 980
$ws.Range("AA8").Value = 70
$ws.Range("AD8").Value = 980
$ws.Range("AE8").Value = 980
$ws.Range("AF8").Value = 14.5
$ws.Range("AG8").Value = 11.5
$ws.Range("AJ8").Value = 980
$ws.Range("AK8").Value = 980
$ws.Range("AM8").Value = 140
$ws.Range("AN8").Value = 980

# Row 9 - Colombian Primera A: Junior FC Barranquilla vs Ind Medellin
$ws.Range("G9").Value = 2.76
$ws.Range("I9").Value = 3.25
$ws.Range("J9").Value = 3.3
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 3.45
$ws.Range("P9").Value = 1.96
$ws.Range("Q9").Value = 1.84
$ws.Range("R9").Value = 1.38
$ws.Range("S9").Value = 2.98
$ws.Range("V9").Value = 1.44
$ws.Range("W9").Value = 1.56
$ws.Range("AB9").Value = 14
$ws.Range("AC9").Value = 1000
